# Auto-generated edit script: updates cryptos list (Price + Volume(1h) columns)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.499.47"
$ws.Range("E2").Value = "  -0.28%  "
$ws.Range("D3").Value = "2.640.49"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.65"
$ws.Range("E5").Value = "  -0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.36"
$ws.Range("E6").Value = "  -1.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.627"
$ws.Range("E8").Value = "  +2.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.126"
$ws.Range("E9").Value = "  +3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.80"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.394"
$ws.Range("E11").Value = "  -1.17%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.58"
$ws.Range("E13").Value = "  -3.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000195"
$ws.Range("E14").Value = "  -0.40%  "
$ws.Range("D15").Value = "3.116.82"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").Value = "65.337.46"
$ws.Range("E16").Value = "  -0.23%  "
$ws.Range("D17").Value = "2.647.36"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.52"
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.70"
$ws.Range("E19").Value = "  -2.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.40"
$ws.Range("E20").Value = "  -2.30%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "347.01"
$ws.Range("E21").Value = "  -1.33%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.82"
$ws.Range("E23").Value = "  -1.20%  "
$ws.Range("E24").Value = "  +0.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.58"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.68"
$ws.Range("E26").Value = "  +2.28%  "
$ws.Range("E27").Value = "  -1.66%  "
$ws.Range("E28").Value = "  -2.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.85"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.11"
$ws.Range("E31").Value = "  -1.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "526.34"
$ws.Range("E32").Value = "  -2.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.74"
$ws.Range("E33").Value = "  -1.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.36"
$ws.Range("E34").Value = "  -3.26%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.37"
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.418"
$ws.Range("E36").Value = "  -1.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.26"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.91"
$ws.Range("E39").Value = "  -1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.91"
$ws.Range("E40").Value = "  -1.93%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "159.73"
$ws.Range("E42").Value = "  -3.50%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.04"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0601"
$ws.Range("E44").Value = "  -2.34%  "
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.49"
$ws.Range("E46").Value = "  -3.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.632"
$ws.Range("E47").Value = "  -2.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0253"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0992"
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").Value = "0.0₆0250"
$ws.Range("E50").Value = "  +7.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.61"
$ws.Range("E51").Value = "  -2.13%  "
